$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently follows
#    the title (H1) paragraph.
# ------------------------------------------------------------------
$metaParaText = 'Meta description: Discover the pros and cons of playing Eye of Horus Megaways, an Ancient Egypt-themed slot game with free spins and chances to win big.'
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pText = $p.Range.Text.TrimEnd([char]13)
    if ($pText -eq $metaParaText -or $pText.StartsWith('Meta description')) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Eye of Horus Megaways Free |
#    Pros, Cons, & Gameplay Review") right before the very last
#    paragraph in the document (the former "Prompt:" paragraph).
#    We clone an existing plain-style paragraph (including its
#    paragraph mark) so the new paragraph picks up regular "Normal"
#    formatting instead of inheriting the italic run that follows it.
# ------------------------------------------------------------------
$template = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Style.NameLocal -eq "Normal" -and $candidate.Range.Text.TrimEnd([char]13) -ne "") {
        $template = $candidate.Range
        break
    }
}
$clonedText = $template.FormattedText

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPoint = $lastPara.Range.Start
$insertTarget = $d.Range($insertPoint, $insertPoint)
$insertTarget.FormattedText = $clonedText

$newParaIndex = $count          # the cloned paragraph is now here, pushing the old content down by one
$newPara = $d.Paragraphs.Item($newParaIndex)
$newTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newTextRange.Text = "Play Eye of Horus Megaways Free | Pros, Cons, & Gameplay Review"

$newBoldRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newBoldRange.Font.Bold = 1

# ------------------------------------------------------------------
# 3) Replace the old "Prompt: ..." image-generation text (now the
#    last paragraph) with the real meta-description copy, keeping
#    the existing italic formatting.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$oldText = 'Prompt: Create a feature image for Eye of Horus Megaways that captures the excitement and adventure of Ancient Egypt, while also including a fun and quirky twist. The image should be in a cartoon style and should prominently feature a happy Maya warrior wearing glasses. The warrior can be holding a staff or any other symbol that is commonly associated with Ancient Egypt, and should be set against a backdrop of pyramids and hieroglyphics. The image should be colorful and eye-catching, with a sense of humor that will appeal to players who are looking for a fun and entertaining slot machine experience. The text "Eye of Horus Megaways" should be included somewhere on the image, along with a tagline that conveys the excitement and adventure of playing this game.'
$newText = 'Discover the pros and cons of playing Eye of Horus Megaways, an Ancient Egypt-themed slot game with free spins and chances to win big.'

$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
